$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the data range to Text format first so numeric-looking strings
# (e.g. "233.56", "1.00") are stored as text, matching the source file
# (every data cell there is stored as text), instead of being
# auto-converted to numbers by the input parser.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '37.375.19'
$ws.Range("E2").Value = '  +2.27%  '
$ws.Range("D3").Value = '2.058.10'
$ws.Range("E3").Value = '  +1.78%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '233.56'
$ws.Range("E5").Value = '  -0.63%  '
$ws.Range("E6").Value = '  +2.70%  '
$ws.Range("D7").Value = '57.90'
$ws.Range("E7").Value = '  +5.50%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").Value = '0.384'
$ws.Range("E9").Value = '  +3.67%  '
$ws.Range("D10").Value = '58.75'
$ws.Range("E10").Value = '  +2.11%  '
$ws.Range("D11").Value = '0.0763'
$ws.Range("E11").Value = '  +1.78%  '
$ws.Range("E12").Value = '  +0.36%  '
$ws.Range("D13").Value = '2.359.98'
$ws.Range("E13").Value = '  +1.65%  '
$ws.Range("D14").Value = '14.37'
$ws.Range("E14").Value = '  +1.23%  '
$ws.Range("D15").Value = '20.88'
$ws.Range("E15").Value = '  +3.79%  '
$ws.Range("E16").Value = '  +1.27%  '
$ws.Range("D17").Value = '5.17'
$ws.Range("E17").Value = '  +1.43%  '
$ws.Range("D18").Value = '2.036.54'
$ws.Range("E18").Value = '  +0.78%  '
$ws.Range("D19").Value = '37.576.71'
$ws.Range("E19").Value = '  +3.05%  '
$ws.Range("D20").Value = '6.22'
$ws.Range("E20").Value = '  +15.34%  '
$ws.Range("D21").Value = '69.25'
$ws.Range("E21").Value = '  +2.28%  '
$ws.Range("D22").Value = '0.0₃0812'
$ws.Range("E22").Value = '  +1.88%  '
$ws.Range("D23").Value = '226.12'
$ws.Range("E23").Value = '  +2.60%  '
$ws.Range("E24").Value = '  +0.04%  '
$ws.Range("B25").Value = 'PancakeSwap'
$ws.Range("C25").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D25").Value = '2.43'
$ws.Range("E25").Value = '  +1.12%  '
$ws.Range("B26").Value = 'Toncoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D26").Value = '2.39'
$ws.Range("D27").Value = '165.64'
$ws.Range("E27").Value = '  +1.50%  '
$ws.Range("D28").Value = '1.47'
$ws.Range("E28").Value = '  +5.79%  '
$ws.Range("D29").Value = '8.93'
$ws.Range("E29").Value = '  +3.83%  '
$ws.Range("E30").Value = '  +1.12%  '
$ws.Range("D31").Value = '19.15'
$ws.Range("E31").Value = '  +1.12%  '
$ws.Range("D32").Value = '0.119'
$ws.Range("E32").Value = '  +1.39%  '
$ws.Range("D33").Value = '4.49'
$ws.Range("E33").Value = '  +2.91%  '
$ws.Range("B34").Value = 'LidoDAOToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D34").Value = '2.58'
$ws.Range("E34").Value = '  +5.52%  '
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").Value = '0.0621'
$ws.Range("E35").Value = '  +2.81%  '
$ws.Range("E36").Value = '  +8.17%  '
$ws.Range("D37").Value = '1.00'
$ws.Range("D38").Value = '5.95'
$ws.Range("E38").Value = '  +4.55%  '
$ws.Range("D39").Value = '3.32'
$ws.Range("E39").Value = '  +0.31%  '
$ws.Range("E40").Value = '  -0.85%  '
$ws.Range("D41").Value = '4.76'
$ws.Range("E41").Value = '  +14.16%  '
$ws.Range("D43").Value = '0.0946'
$ws.Range("E43").Value = '  +1.94%  '
$ws.Range("D44").Value = '96.02'
$ws.Range("E44").Value = '  +6.59%  '
$ws.Range("D45").Value = '1.458.41'
$ws.Range("E45").Value = '  +0.16%  '
$ws.Range("D46").Value = '1.18'
$ws.Range("E46").Value = '  +6.02%  '
$ws.Range("E47").Value = '  +4.02%  '
$ws.Range("D48").Value = '15.68'
$ws.Range("E48").Value = '  +2.15%  '
$ws.Range("E49").Value = '  +1.70%  '
$ws.Range("D50").Value = '7.19'
$ws.Range("E50").Value = '  +4.80%  '
$ws.Range("E51").Value = '  +1.99%  '

# Drop the temporary Text number format again so the cell styling
# matches the original workbook (plain default style everywhere).
$dataRange.ClearFormats()
